# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.790.21"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "1.664.59"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'330.43"
$ws.Range("E5").Value = "  +8.28%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.3652"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("D8").Value = "'47.32"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.3238"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "'1.140"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").Value = "'0.07066"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'6.067"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").Value = "'19.60"
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "1.663.64"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "'6.606"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "'0.00001052"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "'0.06605"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'78.69"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "'5.937"
$ws.Range("E21").Value = "  +1.06%  "
$ws.Range("D22").Value = "'15.82"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "'12.49"
$ws.Range("E23").Value = "  +3.57%  "
$ws.Range("D24").Value = "24.782.73"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'2.450"
$ws.Range("E25").Value = "  +2.50%  "
$ws.Range("D26").Value = "'2.421"
$ws.Range("E26").Value = "  +4.10%  "
$ws.Range("D27").Value = "'148.47"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").Value = "'18.66"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "1.847.01"
$ws.Range("D30").Value = "'125.69"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "'1.173"
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("D32").Value = "'4.074"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").Value = "'5.708"
$ws.Range("E33").Value = "  +2.52%  "
$ws.Range("D34").Value = "'0.08495"
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("D35").Value = "'1.634"
$ws.Range("E35").Value = "  -1.74%  "
$ws.Range("D36").Value = "'12.15"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").Value = "'5.161"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'0.02267"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").Value = "'0.06075"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2083"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.218"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "'8.229"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.26%  "
$ws.Range("D44").Value = "'0.5924"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("D45").Value = "'13.50"
$ws.Range("E45").Value = "  +7.49%  "
$ws.Range("D46").Value = "'3.854"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("D47").Value = "'0.5666"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").Value = "'125.63"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("D49").Value = "'1.948"
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("D50").Value = "'0.06969"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'1.189"
$ws.Range("E51").Value = "  +4.44%  "
